$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new rows before current row 4 ("options"), shifting existing
# rows 4:33 down to 7:36.
$ws.Rows.Item(4).Resize(3).Insert()

# Fill the 3 new rows with "play"/"PLAY", "credits"/"CREDITS",
# "credits_detail"/long text, matching the style pattern of column B
# (style index 2 = wrapText) used by the other header-like rows.
$ws.Range("A4").Value = "play"
$ws.Range("B4").Value = "PLAY"

$ws.Range("A5").Value = "credits"
$ws.Range("B5").Value = "CREDITS"

$ws.Range("A6").Value = "credits_detail"
$ws.Range("B6").Value = "Written By: David Dionisio\nMusic From: Kevin Macleod"

# Apply the same style (wrap-text, as used on B3/B4.. for short labels)
# to the B cells of the new rows by copying the style from B7 (the row
# that used to be B4, "OPTIONS", which already carries the style).
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B4:B6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the sheet view to match: select B6, scroll to top-left A1.
$ws.Range("B6").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
